$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(8, 4).Range.Text = "92-29=63"
$tbl.Cell(18, 5).Range.Text = "12-3=9"

$d.Content.Find.Execute("64-56=8", $true, $false, $false, $false, $false, $true, 1, $false, "6+56=62", 2)
$d.Content.Find.Execute("49+29=78", $true, $false, $false, $false, $false, $true, 1, $false, "18-13=5", 2)
$d.Content.Find.Execute("0+29=29", $true, $false, $false, $false, $false, $true, 1, $false, "38+22=60", 2)
$d.Content.Find.Execute("49+23=72", $true, $false, $false, $false, $false, $true, 1, $false, "71-29=42", 2)
$d.Content.Find.Execute("47-44=3", $true, $false, $false, $false, $false, $true, 1, $false, "85-62=23", 2)
$d.Content.Find.Execute("28-8=20", $true, $false, $false, $false, $false, $true, 1, $false, "83-7=76", 2)
$d.Content.Find.Execute("14+79=93", $true, $false, $false, $false, $false, $true, 1, $false, "74-0=74", 2)
$d.Content.Find.Execute("75-48=27", $true, $false, $false, $false, $false, $true, 1, $false, "28+29=57", 2)
$d.Content.Find.Execute("75-3=72", $true, $false, $false, $false, $false, $true, 1, $false, "18+42=60", 2)
$d.Content.Find.Execute("30-11=19", $true, $false, $false, $false, $false, $true, 1, $false, "6+82=88", 2)
$d.Content.Find.Execute("4+86=90", $true, $false, $false, $false, $false, $true, 1, $false, "72-45=27", 2)
$d.Content.Find.Execute("47+21=68", $true, $false, $false, $false, $false, $true, 1, $false, "12+46=58", 2)
$d.Content.Find.Execute("93-49=44", $true, $false, $false, $false, $false, $true, 1, $false, "50-21=29", 2)
$d.Content.Find.Execute("58-2=56", $true, $false, $false, $false, $false, $true, 1, $false, "15+38=53", 2)
$d.Content.Find.Execute("16+8=24", $true, $false, $false, $false, $false, $true, 1, $false, "82-60=22", 2)
$d.Content.Find.Execute("88-66=22", $true, $false, $false, $false, $false, $true, 1, $false, "57+30=87", 2)
$d.Content.Find.Execute("40+3=43", $true, $false, $false, $false, $false, $true, 1, $false, "57-12=45", 2)
$d.Content.Find.Execute("57+9=66", $true, $false, $false, $false, $false, $true, 1, $false, "37-22=15", 2)
$d.Content.Find.Execute("5+89=94", $true, $false, $false, $false, $false, $true, 1, $false, "46-22=24", 2)
$d.Content.Find.Execute("80-1=79", $true, $false, $false, $false, $false, $true, 1, $false, "26-10=16", 2)
$d.Content.Find.Execute("89+7=96", $true, $false, $false, $false, $false, $true, 1, $false, "43-43=0", 2)
$d.Content.Find.Execute("77+9=86", $true, $false, $false, $false, $false, $true, 1, $false, "49-25=24", 2)
$d.Content.Find.Execute("83-49=34", $true, $false, $false, $false, $false, $true, 1, $false, "64+5=69", 2)
$d.Content.Find.Execute("18-18=0", $true, $false, $false, $false, $false, $true, 1, $false, "49-8=41", 2)
$d.Content.Find.Execute("75-51=24", $true, $false, $false, $false, $false, $true, 1, $false, "72-13=59", 2)
$d.Content.Find.Execute("92-43=49", $true, $false, $false, $false, $false, $true, 1, $false, "80-23=57", 2)
$d.Content.Find.Execute("68-35=33", $true, $false, $false, $false, $false, $true, 1, $false, "5+56=61", 2)
$d.Content.Find.Execute("57-46=11", $true, $false, $false, $false, $false, $true, 1, $false, "94-60=34", 2)
$d.Content.Find.Execute("15+55=70", $true, $false, $false, $false, $false, $true, 1, $false, "23+14=37", 2)
$d.Content.Find.Execute("23-14=9", $true, $false, $false, $false, $false, $true, 1, $false, "21-6=15", 2)
$d.Content.Find.Execute("58-16=42", $true, $false, $false, $false, $false, $true, 1, $false, "61+19=80", 2)
$d.Content.Find.Execute("54-7=47", $true, $false, $false, $false, $false, $true, 1, $false, "91-8=83", 2)
$d.Content.Find.Execute("27+20=47", $true, $false, $false, $false, $false, $true, 1, $false, "65+5=70", 2)
$d.Content.Find.Execute("94-64=30", $true, $false, $false, $false, $false, $true, 1, $false, "29-26=3", 2)
$d.Content.Find.Execute("74-23=51", $true, $false, $false, $false, $false, $true, 1, $false, "24+13=37", 2)
$d.Content.Find.Execute("65-15=50", $true, $false, $false, $false, $false, $true, 1, $false, "20+42=62", 2)
$d.Content.Find.Execute("46-7=39", $true, $false, $false, $false, $false, $true, 1, $false, "87-40=47", 2)
$d.Content.Find.Execute("98-27=71", $true, $false, $false, $false, $false, $true, 1, $false, "84-13=71", 2)
$d.Content.Find.Execute("44+41=85", $true, $false, $false, $false, $false, $true, 1, $false, "65-10=55", 2)
$d.Content.Find.Execute("90-82=8", $true, $false, $false, $false, $false, $true, 1, $false, "99-27=72", 2)
$d.Content.Find.Execute("55-4=51", $true, $false, $false, $false, $false, $true, 1, $false, "4+81=85", 2)
$d.Content.Find.Execute("60+19=79", $true, $false, $false, $false, $false, $true, 1, $false, "23+20=43", 2)
$d.Content.Find.Execute("19+52=71", $true, $false, $false, $false, $false, $true, 1, $false, "1+91=92", 2)
$d.Content.Find.Execute("2+81=83", $true, $false, $false, $false, $false, $true, 1, $false, "61+11=72", 2)
$d.Content.Find.Execute("51-35=16", $true, $false, $false, $false, $false, $true, 1, $false, "32+0=32", 2)
$d.Content.Find.Execute("65+13=78", $true, $false, $false, $false, $false, $true, 1, $false, "63-34=29", 2)
$d.Content.Find.Execute("68+13=81", $true, $false, $false, $false, $false, $true, 1, $false, "39+48=87", 2)
$d.Content.Find.Execute("86-46=40", $true, $false, $false, $false, $false, $true, 1, $false, "7+64=71", 2)
$d.Content.Find.Execute("86-69=17", $true, $false, $false, $false, $false, $true, 1, $false, "29+52=81", 2)
$d.Content.Find.Execute("71-52=19", $true, $false, $false, $false, $false, $true, 1, $false, "31-30=1", 2)
$d.Content.Find.Execute("67+6=73", $true, $false, $false, $false, $false, $true, 1, $false, "12+82=94", 2)
$d.Content.Find.Execute("24-3=21", $true, $false, $false, $false, $false, $true, 1, $false, "92+5=97", 2)
$d.Content.Find.Execute("13+64=77", $true, $false, $false, $false, $false, $true, 1, $false, "54-34=20", 2)
$d.Content.Find.Execute("11+88=99", $true, $false, $false, $false, $false, $true, 1, $false, "43+23=66", 2)
$d.Content.Find.Execute("48-7=41", $true, $false, $false, $false, $false, $true, 1, $false, "49+18=67", 2)
$d.Content.Find.Execute("57-44=13", $true, $false, $false, $false, $false, $true, 1, $false, "0+95=95", 2)
$d.Content.Find.Execute("11+60=71", $true, $false, $false, $false, $false, $true, 1, $false, "89-56=33", 2)
$d.Content.Find.Execute("65+19=84", $true, $false, $false, $false, $false, $true, 1, $false, "37-13=24", 2)
$d.Content.Find.Execute("85-27=58", $true, $false, $false, $false, $false, $true, 1, $false, "79+17=96", 2)
$d.Content.Find.Execute("26-16=10", $true, $false, $false, $false, $false, $true, 1, $false, "85-59=26", 2)
$d.Content.Find.Execute("72-20=52", $true, $false, $false, $false, $false, $true, 1, $false, "78-60=18", 2)
$d.Content.Find.Execute("99-29=70", $true, $false, $false, $false, $false, $true, 1, $false, "14+76=90", 2)
$d.Content.Find.Execute("49-38=11", $true, $false, $false, $false, $false, $true, 1, $false, "51-34=17", 2)
$d.Content.Find.Execute("28+37=65", $true, $false, $false, $false, $false, $true, 1, $false, "38+32=70", 2)
$d.Content.Find.Execute("76-6=70", $true, $false, $false, $false, $false, $true, 1, $false, "1+58=59", 2)
$d.Content.Find.Execute("31+46=77", $true, $false, $false, $false, $false, $true, 1, $false, "77+3=80", 2)
$d.Content.Find.Execute("49-9=40", $true, $false, $false, $false, $false, $true, 1, $false, "99-74=25", 2)
$d.Content.Find.Execute("57+33=90", $true, $false, $false, $false, $false, $true, 1, $false, "64+27=91", 2)
$d.Content.Find.Execute("97-72=25", $true, $false, $false, $false, $false, $true, 1, $false, "99-57=42", 2)
$d.Content.Find.Execute("28+1=29", $true, $false, $false, $false, $false, $true, 1, $false, "70-70=0", 2)
$d.Content.Find.Execute("19+12=31", $true, $false, $false, $false, $false, $true, 1, $false, "39+13=52", 2)
$d.Content.Find.Execute("30-25=5", $true, $false, $false, $false, $false, $true, 1, $false, "3+55=58", 2)
$d.Content.Find.Execute("80-40=40", $true, $false, $false, $false, $false, $true, 1, $false, "72-61=11", 2)
$d.Content.Find.Execute("34-33=1", $true, $false, $false, $false, $false, $true, 1, $false, "18+49=67", 2)
$d.Content.Find.Execute("12+59=71", $true, $false, $false, $false, $false, $true, 1, $false, "70-10=60", 2)
$d.Content.Find.Execute("94-28=66", $true, $false, $false, $false, $false, $true, 1, $false, "18-17=1", 2)
$d.Content.Find.Execute("66-15=51", $true, $false, $false, $false, $false, $true, 1, $false, "86-13=73", 2)
$d.Content.Find.Execute("37-8=29", $true, $false, $false, $false, $false, $true, 1, $false, "12+8=20", 2)
$d.Content.Find.Execute("93-42=51", $true, $false, $false, $false, $false, $true, 1, $false, "44-32=12", 2)
$d.Content.Find.Execute("47+33=80", $true, $false, $false, $false, $false, $true, 1, $false, "29+49=78", 2)
$d.Content.Find.Execute("6+30=36", $true, $false, $false, $false, $false, $true, 1, $false, "22+22=44", 2)
$d.Content.Find.Execute("84+11=95", $true, $false, $false, $false, $false, $true, 1, $false, "10+44=54", 2)
$d.Content.Find.Execute("1+42=43", $true, $false, $false, $false, $false, $true, 1, $false, "3+73=76", 2)
$d.Content.Find.Execute("86-33=53", $true, $false, $false, $false, $false, $true, 1, $false, "61+26=87", 2)
$d.Content.Find.Execute("33-22=11", $true, $false, $false, $false, $false, $true, 1, $false, "33+34=67", 2)
$d.Content.Find.Execute("16+20=36", $true, $false, $false, $false, $false, $true, 1, $false, "95-25=70", 2)
$d.Content.Find.Execute("54-10=44", $true, $false, $false, $false, $false, $true, 1, $false, "56-16=40", 2)
$d.Content.Find.Execute("39-11=28", $true, $false, $false, $false, $false, $true, 1, $false, "39+5=44", 2)
$d.Content.Find.Execute("46+16=62", $true, $false, $false, $false, $false, $true, 1, $false, "99-98=1", 2)
$d.Content.Find.Execute("40-24=16", $true, $false, $false, $false, $false, $true, 1, $false, "91-59=32", 2)
$d.Content.Find.Execute("3+81=84", $true, $false, $false, $false, $false, $true, 1, $false, "62-28=34", 2)
$d.Content.Find.Execute("30+51=81", $true, $false, $false, $false, $false, $true, 1, $false, "92-60=32", 2)
$d.Content.Find.Execute("40-12=28", $true, $false, $false, $false, $false, $true, 1, $false, "34+40=74", 2)
$d.Content.Find.Execute("5+38=43", $true, $false, $false, $false, $false, $true, 1, $false, "3+66=69", 2)
$d.Content.Find.Execute("1+44=45", $true, $false, $false, $false, $false, $true, 1, $false, "99-31=68", 2)
$d.Content.Find.Execute("36+54=90", $true, $false, $false, $false, $false, $true, 1, $false, "75-59=16", 2)
$d.Content.Find.Execute("65-47=18", $true, $false, $false, $false, $false, $true, 1, $false, "87-76=11", 2)
$d.Content.Find.Execute("15+50=65", $true, $false, $false, $false, $false, $true, 1, $false, "53+14=67", 2)
